$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 362.0625
$ws.Range("I33").Value = 311
$ws.Range("K33").Value = 311
$ws.Range("M33").Value = -82
$ws.Range("H53").Value = 1197.4166
$ws.Range("J53").Value = 656.875
$ws.Range("L53").Value = 656.875
$ws.Range("N53").Value = -1930.875
$ws.Range("H86").Value = 69447870
$ws.Range("I86").Value = 25003254
$ws.Range("J86").Value = 125003640
$ws.Range("K86").Value = 25003254
$ws.Range("L86").Value = 125003640
$ws.Range("M86").Value = -25002131
$ws.Range("N86").Value = -125005886
$ws.Range("H89").Value = 69447870
$ws.Range("I89").Value = 25003254
$ws.Range("J89").Value = 125003640
$ws.Range("K89").Value = 125016270
$ws.Range("L89").Value = 625018200
$ws.Range("M89").Value = -125010654
$ws.Range("N89").Value = -625029432
$ws.Range("H92").Value = 410
$ws.Range("I92").Value = 382
$ws.Range("K92").Value = 382
$ws.Range("M92").Value = 866
$ws.Range("H96").Value = 3723.5
$ws.Range("I96").Value = 462.4
$ws.Range("K96").Value = 1387.2
$ws.Range("M96").Value = -14.19999999999982
$ws.Range("H111").Value = 564
$ws.Range("I111").Value = 470.4
$ws.Range("J111").Value = 1500
$ws.Range("K111").Value = 1411.2
$ws.Range("L111").Value = 4500
$ws.Range("M111").Value = 1655.8
$ws.Range("N111").Value = -10634
$ws.Range("H112").Value = 1161.1786
$ws.Range("I112").Value = 520.8
$ws.Range("J112").Value = 1300.3914
$ws.Range("K112").Value = 1562.4
$ws.Range("L112").Value = 3901.1742
$ws.Range("M112").Value = -454.3999999999999
$ws.Range("N112").Value = -6117.174199999999
$ws.Range("H138").Value = 1920.95
$ws.Range("I138").Value = 1617.1538
$ws.Range("J138").Value = 2485.1428
$ws.Range("K138").Value = 4851.4614
$ws.Range("L138").Value = 7455.428400000001
$ws.Range("M138").Value = 288.5385999999999
$ws.Range("N138").Value = -17735.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2404.8572
$ws.Range("I61").Value = 1416.25
$ws.Range("K61").Value = 1416.25
$ws.Range("M61").Value = -1204.25
$ws.Range("H98").Value = 44850
$ws.Range("J98").Value = 57500
$ws.Range("L98").Value = 57500
$ws.Range("N98").Value = -63490
$ws.Range("H102").Value = 100974.63
$ws.Range("I102").Value = 143721.28
$ws.Range("J102").Value = 26168
$ws.Range("K102").Value = 143721.28
$ws.Range("L102").Value = 26168
$ws.Range("M102").Value = -142099.28
$ws.Range("N102").Value = -29412
$ws.Range("H124").Value = 62419
$ws.Range("J124").Value = 62419
$ws.Range("L124").Value = 62419
$ws.Range("N124").Value = -72239
$ws.Range("H132").Value = 2653.75
$ws.Range("I132").Value = 1932.5
$ws.Range("K132").Value = 5797.5
$ws.Range("M132").Value = -3267.5
$ws.Range("H136").Value = 2404.8572
$ws.Range("I136").Value = 1416.25
$ws.Range("K136").Value = 4248.75
$ws.Range("M136").Value = -1698.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 294116.25
$ws.Range("I20").Value = 335925.84
$ws.Range("J20").Value = 1449
$ws.Range("K20").Value = 335925.84
$ws.Range("L20").Value = 1449
$ws.Range("M20").Value = -335678.84
$ws.Range("N20").Value = -1943
$ws.Range("H82").Value = 24128
$ws.Range("I82").Value = 9257
$ws.Range("J82").Value = 38999
$ws.Range("K82").Value = 9257
$ws.Range("L82").Value = 38999
$ws.Range("M82").Value = -8874
$ws.Range("N82").Value = -39765
$ws.Range("H85").Value = 24128
$ws.Range("I85").Value = 9257
$ws.Range("J85").Value = 38999
$ws.Range("K85").Value = 9257
$ws.Range("L85").Value = 38999
$ws.Range("M85").Value = -7931
$ws.Range("N85").Value = -41651
$ws.Range("H86").Value = 6982.4736
$ws.Range("I86").Value = 5206.1665
$ws.Range("J86").Value = 10027.571
$ws.Range("K86").Value = 5206.1665
$ws.Range("L86").Value = 10027.571
$ws.Range("M86").Value = -4083.1665
$ws.Range("N86").Value = -12273.571
$ws.Range("H89").Value = 6982.4736
$ws.Range("I89").Value = 5206.1665
$ws.Range("J89").Value = 10027.571
$ws.Range("K89").Value = 26030.8325
$ws.Range("L89").Value = 50137.855
$ws.Range("M89").Value = -20414.8325
$ws.Range("N89").Value = -61369.855
$ws.Range("H110").Value = 84922
$ws.Range("J110").Value = 84922
$ws.Range("L110").Value = 84922
$ws.Range("N110").Value = -93102
$ws.Range("H132").Value = 38135.453
$ws.Range("J132").Value = 38135.453
$ws.Range("L132").Value = 38135.453
$ws.Range("N132").Value = -48255.453
$ws.Range("H140").Value = 43499.383
$ws.Range("J140").Value = 43499.383
$ws.Range("L140").Value = 43499.383
$ws.Range("N140").Value = -53859.383

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 6249.1
$ws.Range("I86").Value = 3865.1667
$ws.Range("K86").Value = 3865.1667
$ws.Range("M86").Value = -2742.1667
$ws.Range("H89").Value = 6249.1
$ws.Range("I89").Value = 3865.1667
$ws.Range("K89").Value = 19325.8335
$ws.Range("M89").Value = -13709.8335
$ws.Range("H105").Value = 283577.5
$ws.Range("I105").Value = 557155
$ws.Range("K105").Value = 557155
$ws.Range("M105").Value = -555408
$ws.Range("H134").Value = 3251.6428
$ws.Range("J134").Value = 2987.8
$ws.Range("L134").Value = 8963.400000000001
$ws.Range("N134").Value = -14033.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 6399407
$ws.Range("I4").Value = 5555975
$ws.Range("J4").Value = 11460000
$ws.Range("K4").Value = 16667925
$ws.Range("L4").Value = 34380000
$ws.Range("M4").Value = -16667813
$ws.Range("N4").Value = -34380224
$ws.Range("H48").Value = 1501
$ws.Range("J48").Value = 999
$ws.Range("L48").Value = 2997
$ws.Range("N48").Value = -3497
$ws.Range("H50").Value = 1072
$ws.Range("I50").Value = 1072
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 3216
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = -2735
$ws.Range("N50").ClearContents()
$ws.Range("H53").Value = 1072
$ws.Range("I53").Value = 1072
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 3216
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = -2735
$ws.Range("N53").ClearContents()
$ws.Range("H56").Value = 6032.4
$ws.Range("I56").Value = 6032.4
$ws.Range("K56").Value = 6032.4
$ws.Range("M56").Value = -5502.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("H102").Value = 2298.1
$ws.Range("I102").Value = 2243.7144
$ws.Range("K102").Value = 2243.7144
$ws.Range("M102").Value = -621.7143999999998
$ws.Range("H123").Value = 46557
$ws.Range("J123").Value = 46557
$ws.Range("L123").Value = 46557
$ws.Range("N123").Value = -51457
$ws.Range("H132").Value = 3921.4546
$ws.Range("I132").Value = 3266
$ws.Range("J132").Value = 5068.5
$ws.Range("K132").Value = 9798
$ws.Range("L132").Value = 15205.5
$ws.Range("M132").Value = -7268
$ws.Range("N132").Value = -20265.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 2241.3225
$ws.Range("I55").Value = 1301.9524
$ws.Range("J55").Value = 4214
$ws.Range("K55").Value = 1301.9524
$ws.Range("L55").Value = 4214
$ws.Range("M55").Value = -1128.9524
$ws.Range("N55").Value = -4560
$ws.Range("H82").Value = 2390.625
$ws.Range("I82").Value = 2390.625
$ws.Range("K82").Value = 2390.625
$ws.Range("M82").Value = -2029.625
$ws.Range("H85").Value = 2390.625
$ws.Range("I85").Value = 2390.625
$ws.Range("K85").Value = 2390.625
$ws.Range("M85").Value = -1142.625
$ws.Range("H122").Value = 66670560
$ws.Range("I122").Value = 111114510
$ws.Range("K122").Value = 333343530
$ws.Range("M122").Value = -333341080
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H132").Value = 1414.7742
$ws.Range("I132").Value = 1295.1786
$ws.Range("K132").Value = 3885.5358
$ws.Range("M132").Value = -1355.5358
$ws.Range("H133").Value = 86864
$ws.Range("J133").Value = 86864
$ws.Range("L133").Value = 86864
$ws.Range("N133").Value = -91924

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H132").Value = 1320.7693
$ws.Range("I132").Value = 1264.1666
$ws.Range("K132").Value = 3792.4998
$ws.Range("M132").Value = -1262.4998

